$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad") from 45202 to 45203 for all data rows (C2:C500)
$range = $ws.Range("C2:C500")
foreach ($cell in $range.Cells) {
    if ($cell.Value2 -eq 45202) {
        $cell.Value2 = 45203
    }
}
